$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing data row (65) down into
# the new row (66) for the two columns that carry explicit cell styles:
#   column A -> bold / bordered / centered "Indice" style
#   column E -> custom date-time number format style
# Copying preserves the existing style indexes instead of minting new ones.
$ws.Range("A65").Copy($ws.Range("A66"))
$ws.Range("E65").Copy($ws.Range("E66"))

# Populate the new row's values (row 66 = match index 65)
$ws.Range("A66").Value = 65
$ws.Range("B66").Value = "croatia"
$ws.Range("C66").Value = "hnl"
$ws.Range("D66").Value = "2023-2024"
$ws.Range("E66").Value = 45235.625
$ws.Range("F66").Value = "Varazdin"
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = "D. Zagreb"
$ws.Range("I66").Value = 1
$ws.Range("J66").Value = 5.87
$ws.Range("K66").Value = "29/10/2023 17:13"
$ws.Range("L66").Value = 5.71
$ws.Range("M66").Value = "05/11/2023 14:57"
$ws.Range("N66").Value = 4.14
$ws.Range("O66").Value = "29/10/2023 17:13"
$ws.Range("P66").Value = 4.13
$ws.Range("Q66").Value = "05/11/2023 14:57"
$ws.Range("R66").Value = 1.5
$ws.Range("S66").Value = "29/10/2023 17:13"
$ws.Range("T66").Value = 1.58
$ws.Range("U66").Value = "05/11/2023 14:57"
$ws.Range("V66").Value = "https://www.betexplorer.com/football/croatia/hnl/varazdin-din-zagreb/4pcEZ1ks/"
